$d = $word.ActiveDocument

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$tbl = $d.Tables(1)
$cell = $tbl.Cell(1,2)
$para = $cell.Range.Paragraphs(1).Range
$base = $para.Start
$t = $para.Text
Write-Host "orig:" $t

# locate '0' in '10' (the character we will change) to place bookmark right after it once changed
$monthIdx = $t.IndexOf(". 10. ")
$zeroPos = $base + $monthIdx + 3   # position of '0' char start
Write-Host "zeroPos:" $zeroPos
$chk = $d.Range($zeroPos, $zeroPos+1)
Write-Host "chk (expect 0):" $chk.Text

# Step 1: change '0' -> '1'
$chk.Text = "1"
Write-Host "after step1:" $para.Text
Write-Host "XML after step1:"

# Step 2: insert bookmark right after this position (cursor after typed char)
$bmRange = $d.Range($zeroPos+1, $zeroPos+1)
$d.Bookmarks.Add("_GoBack", $bmRange)
Write-Host "after step2:" $para.Text

# Step 3: change '22' -> '6' (day). Recompute offset fresh.
$t3 = $para.Text
$dayIdx = $t3.IndexOf("22. 11. 2019.")
Write-Host "dayIdx:" $dayIdx
$s3 = $base + $dayIdx
$r3 = $d.Range($s3, $s3+2)
Write-Host "r3 (expect 22):" $r3.Text
$r3.Text = "6"
Write-Host "final:" $para.Text
